$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '36.373.34'
$ws.Cells.Item(2, 5).Value = '  +0.04%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.013.91'
$ws.Cells.Item(3, 5).Value = '  -1.85%  '
$ws.Cells.Item(4, 5).Value = '  -0.11%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '252.25'
$ws.Cells.Item(5, 5).Value = '  +2.83%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.638'
$ws.Cells.Item(6, 5).Value = '  -3.56%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '62.35'
$ws.Cells.Item(7, 5).Value = '  +9.77%  '
$ws.Cells.Item(8, 5).Value = '  -0.07%  '
$ws.Cells.Item(9, 2).Value = 'Cardano'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.370'
$ws.Cells.Item(9, 5).Value = '  +0.44%  '
$ws.Cells.Item(10, 2).Value = 'OKB'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '58.53'
$ws.Cells.Item(10, 5).Value = '  -7.19%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0744'
$ws.Cells.Item(11, 5).Value = '  -0.65%  '
$ws.Cells.Item(12, 5).Value = '  -1.97%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.902'
$ws.Cells.Item(13, 5).Value = '  -2.55%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '14.90'
$ws.Cells.Item(14, 5).Value = '  +3.40%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '2.308.74'
$ws.Cells.Item(15, 5).Value = '  -1.81%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '20.70'
$ws.Cells.Item(16, 5).Value = '  +16.27%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '5.46'
$ws.Cells.Item(17, 5).Value = '  +0.95%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '2.011.66'
$ws.Cells.Item(18, 5).Value = '  -1.73%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '36.319.60'
$ws.Cells.Item(19, 5).Value = '  -0.12%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '72.05'
$ws.Cells.Item(20, 5).Value = '  +0.60%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.0₃0865'
$ws.Cells.Item(21, 5).Value = '  +0.76%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '5.30'
$ws.Cells.Item(22, 5).Value = '  +1.40%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '234.57'
$ws.Cells.Item(23, 5).Value = '  -0.76%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '2.70'
$ws.Cells.Item(24, 5).Value = '  +19.18%  '
$ws.Cells.Item(25, 5).Value = '  -0.16%  '
$ws.Cells.Item(26, 5).Value = '  -1.38%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '9.63'
$ws.Cells.Item(27, 5).Value = '  +2.62%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '163.37'
$ws.Cells.Item(28, 5).Value = '  -0.79%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '19.64'
$ws.Cells.Item(29, 5).Value = '  -1.79%  '
$ws.Cells.Item(30, 5).Value = '  -1.05%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '5.12'
$ws.Cells.Item(31, 5).Value = '  +2.55%  '
$ws.Cells.Item(32, 2).Value = 'ImmutableX'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '1.18'
$ws.Cells.Item(32, 5).Value = '  +0.17%  '
$ws.Cells.Item(33, 2).Value = 'Kaspa'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.109'
$ws.Cells.Item(33, 5).Value = '  +22.88%  '
$ws.Cells.Item(34, 5).Value = '  +3.44%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.0608'
$ws.Cells.Item(35, 5).Value = '  +1.30%  '
$ws.Cells.Item(36, 5).Value = '  +10.36%  '
$ws.Cells.Item(37, 5).Value = '  -0.21%  '
$ws.Cells.Item(38, 5).Value = '  -1.31%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '5.93'
$ws.Cells.Item(39, 5).Value = '  +16.92%  '
$ws.Cells.Item(40, 5).Value = '  +15.51%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.23'
$ws.Cells.Item(41, 5).Value = '  +0.68%  '
$ws.Cells.Item(42, 5).Value = '  +1.09%  '
$ws.Cells.Item(43, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '17.07'
$ws.Cells.Item(43, 5).Value = '  +7.04%  '
$ws.Cells.Item(44, 2).Value = 'VeChain'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.0217'
$ws.Cells.Item(44, 5).Value = '  +0.27%  '
$ws.Cells.Item(45, 2).Value = 'ARBITRUM'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '1.13'
$ws.Cells.Item(45, 5).Value = '  +2.60%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '1.452.99'
$ws.Cells.Item(46, 5).Value = '  +3.26%  '
$ws.Cells.Item(47, 2).Value = 'FraxShare'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '7.91'
$ws.Cells.Item(47, 5).Value = '  +5.75%  '
$ws.Cells.Item(48, 2).Value = 'Aave'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '95.05'
$ws.Cells.Item(48, 5).Value = '  +1.50%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '2.64'
$ws.Cells.Item(49, 5).Value = '  +16.28%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '2.93'
$ws.Cells.Item(50, 5).Value = '  -1.13%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '47.28'
$ws.Cells.Item(51, 5).Value = '  +2.13%  '
